# Fixed more winner odds BUGs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the comment-separator typo ("  //  " -> "  // ") in the M-column
#    formulas. M2 holds its own (non-shared) formula. M3:M11 is a shared
#    formula group anchored at M3 - set the whole range at once so the
#    engine re-establishes the M3:M11 shared group instead of exploding it
#    into per-cell formulas.
$ws.Range("M2").Formula = '=""""&A2&""":  {[TEAM."&E2&".id]: "&H2&", [TEAM."&F2&".id]: "&I2&"},  // "&J2&"-"&K2'
$ws.Range("M3:M11").Formula = '=""""&A3&""":  {[TEAM."&E3&".id]: "&H3&", [TEAM."&F3&".id]: "&I3&"},  // "&J3&"-"&K3'

# 2) The N:S columns in rows 4-11 contained stray, empty, styled cells
#    left over from an earlier fill-down. Remove them entirely.
$ws.Range("N4:S11").Clear()

# 3) Restore the selection to the M2:M11 block (with M2 active).
$ws.Range("M2:M11").Select()
